$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 and J1 (use same style as existing headers, e.g. style of H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I and J, rows 2-26
$values = @{
    2  = @(6, 6)
    3  = @(6, 6)
    4  = @(4, 5)
    5  = @(7, 7)
    6  = @(5, 6)
    7  = @(8, 8)
    8  = @(6, 6)
    9  = @(1, 1)
    10 = @(7, 7)
    11 = @(8, 8)
    12 = @(6, 6)
    13 = @(8, 8)
    14 = @(5, 6)
    15 = @(7, 7)
    16 = @(8, 8)
    17 = @(8, 9)
    18 = @(7, 7)
    19 = @(7, 7)
    20 = @(7, 7)
    21 = @(7, 7)
    22 = @(7, 7)
    23 = @(7, 7)
    24 = @(4, 4)
    25 = @(7, 7)
    26 = @(4, 4)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
